# Update all PSSM score values (B2:K21) on Sheet1 with the refreshed
# values from the "supplemental figures" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'double[,]' 20,10
$data[0,0] = -19.52716579668945
$data[0,1] = 2.438076318809256
$data[0,2] = -19.52716579668945
$data[0,3] = -19.52716579668945
$data[0,4] = -19.52716579668945
$data[0,5] = -19.52716579668945
$data[0,6] = -19.52716579668945
$data[0,7] = -19.52716579668945
$data[0,8] = -19.52716579668945
$data[0,9] = -19.52716579668945
$data[1,0] = -19.52716579668945
$data[1,1] = -19.52716579668945
$data[1,2] = -19.52716579668945
$data[1,3] = -19.52716579668945
$data[1,4] = -19.52716579668945
$data[1,5] = -19.52716579668945
$data[1,6] = -19.52716579668945
$data[1,7] = 2.310254638683485
$data[1,8] = -19.52716579668945
$data[1,9] = -19.52716579668945
$data[2,0] = -19.52716579668945
$data[2,1] = 2.134099533508155
$data[2,2] = 2.869416251983346
$data[2,3] = -19.52716579668945
$data[2,4] = 2.634034067857902
$data[2,5] = -19.52716579668945
$data[2,6] = 1.72300364145387
$data[2,7] = -19.52716579668945
$data[2,8] = 2.08163631855632
$data[2,9] = -19.52716579668945
$data[3,0] = -19.52716579668945
$data[3,1] = 1.041058479659213
$data[3,2] = -19.52716579668945
$data[3,3] = -19.52716579668945
$data[3,4] = -19.52716579668945
$data[3,5] = 2.472330894290328
$data[3,6] = -19.52716579668945
$data[3,7] = -19.52716579668945
$data[3,8] = -19.52716579668945
$data[3,9] = -19.52716579668945
$data[4,0] = -19.52716579668945
$data[4,1] = -19.52716579668945
$data[4,2] = -19.52716579668945
$data[4,3] = -19.52716579668945
$data[4,4] = -19.52716579668945
$data[4,5] = -19.52716579668945
$data[4,6] = -19.52716579668945
$data[4,7] = -19.52716579668945
$data[4,8] = -19.52716579668945
$data[4,9] = -19.52716579668945
$data[5,0] = 2.969142118531138
$data[5,1] = -19.52716579668945
$data[5,2] = -19.52716579668945
$data[5,3] = -19.52716579668945
$data[5,4] = -19.52716579668945
$data[5,5] = -19.52716579668945
$data[5,6] = -19.52716579668945
$data[5,7] = -19.52716579668945
$data[5,8] = -19.52716579668945
$data[5,9] = -19.52716579668945
$data[6,0] = -19.52716579668945
$data[6,1] = -19.52716579668945
$data[6,2] = -19.52716579668945
$data[6,3] = 2.867978411755925
$data[6,4] = -19.52716579668945
$data[6,5] = -19.52716579668945
$data[6,6] = -19.52716579668945
$data[6,7] = -19.52716579668945
$data[6,8] = -19.52716579668945
$data[6,9] = -19.52716579668945
$data[7,0] = 3.605170913899634
$data[7,1] = -19.52716579668945
$data[7,2] = -19.52716579668945
$data[7,3] = -19.52716579668945
$data[7,4] = -19.52716579668945
$data[7,5] = -19.52716579668945
$data[7,6] = -19.52716579668945
$data[7,7] = -19.52716579668945
$data[7,8] = -19.52716579668945
$data[7,9] = -19.52716579668945
$data[8,0] = -19.52716579668945
$data[8,1] = -19.52716579668945
$data[8,2] = -19.52716579668945
$data[8,3] = -19.52716579668945
$data[8,4] = -19.52716579668945
$data[8,5] = -19.52716579668945
$data[8,6] = -19.52716579668945
$data[8,7] = 1.581096385139967
$data[8,8] = -19.52716579668945
$data[8,9] = -19.52716579668945
$data[9,0] = -19.52716579668945
$data[9,1] = -19.52716579668945
$data[9,2] = -19.52716579668945
$data[9,3] = 2.054349726080715
$data[9,4] = -19.52716579668945
$data[9,5] = 2.451040853807902
$data[9,6] = -19.52716579668945
$data[9,7] = -19.52716579668945
$data[9,8] = -19.52716579668945
$data[9,9] = -19.52716579668945
$data[10,0] = -19.52716579668945
$data[10,1] = -19.52716579668945
$data[10,2] = -19.52716579668945
$data[10,3] = -19.52716579668945
$data[10,4] = -19.52716579668945
$data[10,5] = -19.52716579668945
$data[10,6] = -19.52716579668945
$data[10,7] = -19.52716579668945
$data[10,8] = -19.52716579668945
$data[10,9] = -19.52716579668945
$data[11,0] = -19.52716579668945
$data[11,1] = -19.52716579668945
$data[11,2] = -19.52716579668945
$data[11,3] = 1.715045943720436
$data[11,4] = -19.52716579668945
$data[11,5] = -19.52716579668945
$data[11,6] = -19.52716579668945
$data[11,7] = -19.52716579668945
$data[11,8] = 2.161503273819974
$data[11,9] = -19.52716579668945
$data[12,0] = -19.52716579668945
$data[12,1] = -19.52716579668945
$data[12,2] = 1.686546594629527
$data[12,3] = -19.52716579668945
$data[12,4] = -19.52716579668945
$data[12,5] = -19.52716579668945
$data[12,6] = -19.52716579668945
$data[12,7] = -19.52716579668945
$data[12,8] = -19.52716579668945
$data[12,9] = -19.52716579668945
$data[13,0] = -19.52716579668945
$data[13,1] = -19.52716579668945
$data[13,2] = -0.1815244922296166
$data[13,3] = -19.52716579668945
$data[13,4] = -19.52716579668945
$data[13,5] = -19.52716579668945
$data[13,6] = -19.52716579668945
$data[13,7] = -19.52716579668945
$data[13,8] = -19.52716579668945
$data[13,9] = -19.52716579668945
$data[14,0] = -19.52716579668945
$data[14,1] = -19.52716579668945
$data[14,2] = -19.52716579668945
$data[14,3] = -19.52716579668945
$data[14,4] = -19.52716579668945
$data[14,5] = -19.52716579668945
$data[14,6] = -19.52716579668945
$data[14,7] = -19.52716579668945
$data[14,8] = 2.350070774115064
$data[14,9] = -19.52716579668945
$data[15,0] = -19.52716579668945
$data[15,1] = 0.8534104686945665
$data[15,2] = 0.1051638430479449
$data[15,3] = -19.52716579668945
$data[15,4] = -19.52716579668945
$data[15,5] = -19.52716579668945
$data[15,6] = 0.8519301317302455
$data[15,7] = 1.217534207244994
$data[15,8] = 1.559893618027613
$data[15,9] = -19.52716579668945
$data[16,0] = -19.52716579668945
$data[16,1] = -19.52716579668945
$data[16,2] = -19.52716579668945
$data[16,3] = -19.52716579668945
$data[16,4] = -19.52716579668945
$data[16,5] = -19.52716579668945
$data[16,6] = 0.6441930566037184
$data[16,7] = 1.197192338268892
$data[16,8] = 1.698880016773094
$data[16,9] = -19.52716579668945
$data[17,0] = -19.52716579668945
$data[17,1] = -19.52716579668945
$data[17,2] = 1.722957367876367
$data[17,3] = -19.52716579668945
$data[17,4] = -19.52716579668945
$data[17,5] = -19.52716579668945
$data[17,6] = 1.81664288095891
$data[17,7] = 2.056848424077369
$data[17,8] = -19.52716579668945
$data[17,9] = -19.52716579668945
$data[18,0] = -19.52716579668945
$data[18,1] = 1.595624922286588
$data[18,2] = 2.07550634009679
$data[18,3] = -19.52716579668945
$data[18,4] = 3.785799811738229
$data[18,5] = -19.52716579668945
$data[18,6] = 2.198124314665307
$data[18,7] = 1.708902716747838
$data[18,8] = -19.52716579668945
$data[18,9] = 4.32192628089363
$data[19,0] = -19.52716579668945
$data[19,1] = 1.724391159077017
$data[19,2] = -19.52716579668945
$data[19,3] = 2.395866025395081
$data[19,4] = -19.52716579668945
$data[19,5] = 3.167113226702429
$data[19,6] = 2.383923681426174
$data[19,7] = -19.52716579668945
$data[19,8] = -19.52716579668945
$data[19,9] = -19.52716579668945

$ws.Range("B2:K21").Value = $data
